$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ZIP code sample value (row 6 scenario) from 560071 to 110010
$ws.Range("D37").Value = 110010

# Update the email-id sample values (row 7 scenario) replacing the "rohit93m" /
# "rohit" handles with "abc123" / "abc"
$ws.Range("D46").Value = "abc123@gmail.com"
$ws.Range("D47").Value = "abc123@yahoo.com"
$ws.Range("D48").Value = "abc123@rediff.com"
$ws.Range("D50").Value = "abc@gmail.com"
$ws.Range("D53").Value = "abc123@@gmail.com"
$ws.Range("D54").Value = "abc123@gmail"

# Update the view's active selection to match the saved UI state
$ws.Range("F54").Select()
